$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new data rows above the existing row 809 (weekly update: a new
# week's worth of "Cilantro" price observations), pushing the old rows
# 809:873 down to 812:876.
$ws.Rows("809:811").Insert()

# --- New row 809 ---------------------------------------------------------
$ws.Cells.Item(809, 1).Value  = 6
$ws.Cells.Item(809, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(809, 3).Value  = "Metropolitana"
$ws.Cells.Item(809, 4).Value  = 44746
$ws.Cells.Item(809, 5).Value  = 13
$ws.Cells.Item(809, 6).Value  = 100112040
$ws.Cells.Item(809, 7).Value  = "Cilantro"
$ws.Cells.Item(809, 8).Value  = "Sin especificar"
$ws.Cells.Item(809, 9).Value  = "Primera"
$ws.Cells.Item(809, 10).Value = 280
$ws.Cells.Item(809, 11).Value = 9500
$ws.Cells.Item(809, 12).Value = 10000
$ws.Cells.Item(809, 13).Value = 9732
$ws.Cells.Item(809, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(809, 15).Value = "Región Metropolitana"
$ws.Cells.Item(809, 16).Value = 270
$ws.Cells.Item(809, 17).Value = 36
$ws.Cells.Item(809, 18).Value = "Hortaliza"

# --- New row 810 ---------------------------------------------------------
$ws.Cells.Item(810, 1).Value  = 6
$ws.Cells.Item(810, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(810, 3).Value  = "Metropolitana"
$ws.Cells.Item(810, 4).Value  = 44746
$ws.Cells.Item(810, 5).Value  = 13
$ws.Cells.Item(810, 6).Value  = 100112040
$ws.Cells.Item(810, 7).Value  = "Cilantro"
$ws.Cells.Item(810, 8).Value  = "Sin especificar"
$ws.Cells.Item(810, 9).Value  = "Primera"
$ws.Cells.Item(810, 10).Value = 220
$ws.Cells.Item(810, 11).Value = 17000
$ws.Cells.Item(810, 12).Value = 18000
$ws.Cells.Item(810, 13).Value = 17409
$ws.Cells.Item(810, 14).Value = "`$/docena de atados"
$ws.Cells.Item(810, 15).Value = "Región Metropolitana"
$ws.Cells.Item(810, 16).Value = 5803
$ws.Cells.Item(810, 17).Value = 3
$ws.Cells.Item(810, 18).Value = "Hortaliza"

# --- New row 811 ---------------------------------------------------------
$ws.Cells.Item(811, 1).Value  = 6
$ws.Cells.Item(811, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(811, 3).Value  = "Metropolitana"
$ws.Cells.Item(811, 4).Value  = 44746
$ws.Cells.Item(811, 5).Value  = 13
$ws.Cells.Item(811, 6).Value  = 100112040
$ws.Cells.Item(811, 7).Value  = "Cilantro"
$ws.Cells.Item(811, 8).Value  = "Sin especificar"
$ws.Cells.Item(811, 9).Value  = "Segunda"
$ws.Cells.Item(811, 10).Value = 60
$ws.Cells.Item(811, 11).Value = 15000
$ws.Cells.Item(811, 12).Value = 15000
$ws.Cells.Item(811, 13).Value = 15000
$ws.Cells.Item(811, 14).Value = "`$/docena de atados"
$ws.Cells.Item(811, 15).Value = "Región Metropolitana"
$ws.Cells.Item(811, 16).Value = 5000
$ws.Cells.Item(811, 17).Value = 3
$ws.Cells.Item(811, 18).Value = "Hortaliza"
